$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")
$rng = $ws.Range("A2:A12")
$rng.Font.Bold = $false
$rng.Font.Bold = $true
